$wb = $excel.ActiveWorkbook

# --- 1. Update the shared "Ready for handoff" status text to "In Translation" ---
# This value shows up on the Overview sheet (columns E2/F2) and on each of the
# per-locale sheets (zh-cn, de-de) in their "Status" column (C2).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the left of -eq; Excel returns
        # TRUE/FALSE cells as real booleans via Value2, and PowerShell's -eq
        # coerces the right-hand side to the left-hand side's type, so a
        # boolean $true on the left would turn "Ready for handoff" -eq into
        # a (wrong) boolean comparison.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value2 = "In Translation"
        }
    }
}

# --- 2. Narrow the status columns ---
# Overview sheet: columns E and F (the per-locale status columns)
# zh-cn / de-de sheets: column C (the "Status" column)
$newWidth = 12.55

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = $newWidth
